$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.821.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.77%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.886.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.71%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.53%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'334.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.71%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.46%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4726"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.66%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3935"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'47.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.60%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.08083"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.77%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'1.026"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.41%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'22.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.14%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.878.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.53%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.996"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.28%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.146"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.02%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.010"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.60%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.06750"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.67%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.00001052"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.00%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'87.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.40%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.14%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.007"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.44%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'27.845.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.83%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.525"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.13%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.12%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.324"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.26%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.110.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.26%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'159.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.85%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'20.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.20%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.108"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.03%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.595"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.61%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'122.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.36%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.9821"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.46%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.09488"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.78%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.452"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.70%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.622"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.04%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.361"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.82%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.06163"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.08%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.84%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.220"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.27%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'8.111"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.15%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +1.33%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.53%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +1.35%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.259"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.87%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.5719"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.63%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'12.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.91%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.949"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.68%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.10%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.06914"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.35%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'113.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.44%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.00000000304"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +7.37%  "
$ws.Range("E51").Style = "Normal"
